$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; everything currently at row 14 downward
# shifts down by one (old row 14 -> 15, ..., old row 60 -> 61).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly price record.
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C14").Value = "Los Lagos"
$ws.Range("D14").Value = 45051
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 100112012
$ws.Range("G14").Value = "Espinaca"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 13000
$ws.Range("N14").Value = "`$/cuna 10 kilos"
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 1300
$ws.Range("Q14").Value = 10
$ws.Range("R14").Value = "Hortaliza"
